$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the new "PLANES_INFO" sheet right after the existing "PLANES" sheet
# ---------------------------------------------------------------------------
$planes = $wb.Worksheets.Item("PLANES")
$infoSheet = $wb.Worksheets.Add($null, $planes)
$infoSheet.Name = "PLANES_INFO"

# Header row (new shared strings: PLAN, PRECIO_USD, DURACION_DIAS, DESCRIPCION CORTA)
$infoSheet.Range("A1").Value = "PLAN"
$infoSheet.Range("B1").Value = "PRECIO_USD"
$infoSheet.Range("C1").Value = "DURACION_DIAS"
$infoSheet.Range("D1").Value = "DESCRIPCION CORTA"

# Row 2 - Basico
$infoSheet.Range("A2").Value = "Basico"
$infoSheet.Range("B2").Value = 4.99
$infoSheet.Range("C2").Value = 30
$infoSheet.Range("D2").Value = "Análisis básico con indicadores clave"

# Row 3 - Pro
$infoSheet.Range("A3").Value = "Pro"
$infoSheet.Range("B3").Value = 9.99
$infoSheet.Range("C3").Value = 30
$infoSheet.Range("D3").Value = "Todo el análisis avanzado + exportaciones"

# Row 4 - Premium
$infoSheet.Range("A4").Value = "Premium"
$infoSheet.Range("B4").Value = 19.99
$infoSheet.Range("C4").Value = 30
$infoSheet.Range("D4").Value = "IA + gráficos inteligentes + comparativos"

# Match the formatting used elsewhere in the workbook: header row styled like
# EN_ANALISIS!A1 (bold, centered, wrapped) and data rows styled like
# GRAFICAS!A5 (centered vertically, wrapped).
$enAnalisis = $wb.Worksheets.Item("EN_ANALISIS")
$graficas = $wb.Worksheets.Item("GRAFICAS")

$enAnalisis.Range("A1").Copy()
$infoSheet.Range("A1:D1").PasteSpecial(-4122)

$graficas.Range("A5").Copy()
$infoSheet.Range("A2:D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths (closest values reachable through the ColumnWidth property,
# which is quantized to the workbook's character-width grid).
$infoSheet.Columns.Item(1).ColumnWidth = 7.333333333333333
$infoSheet.Columns.Item(2).ColumnWidth = 11
$infoSheet.Columns.Item(3).ColumnWidth = 14.333333333333332
$infoSheet.Columns.Item(4).ColumnWidth = 34.666666666666664

# ---------------------------------------------------------------------------
# 2) Clean up the now-redundant explicit styles on the "PLANES" sheet
#    (B2,B3,B4,B6 -> same look as GRAFICAS!A5 ; B5,B7,B8,B9,B10 -> default)
# ---------------------------------------------------------------------------
$graficas.Range("A5").Copy()
foreach ($addr in @("B2", "B3", "B4", "B6")) {
    $planes.Range($addr).PasteSpecial(-4122)
}

$planes.Range("C2").Copy()
foreach ($addr in @("B5", "B7", "B8", "B9", "B10")) {
    $planes.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Make the new sheet the active / selected tab (matches the workbook's
#    bookViews activeTab + the new sheet's tabSelected flag)
# ---------------------------------------------------------------------------
$infoSheet.Activate()
$infoSheet.Range("A1").Select()
